$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B8").Value = "2025-01-15T07:47:50+00:00"

$ws.Range("B15").Formula = "=""true"""
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
